# NYPD CompStat 84th Precinct weekly update
# - Bumps the report volume/week header strings to the next week.
# - Refreshes the crime-stat grid (rows 15-30) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text (volume number / reporting week) ----
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# ---- Row 15 (Rape) ----
$ws.Range("N15").Value = -81.818181818181

# ---- Row 16 (Robbery) ----
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 250
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 95
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = -1.041666666666
$ws.Range("L16").Value = 15.853658536585
$ws.Range("M16").Value = -18.803418803418
$ws.Range("N16").Value = -88.18407960199

# ---- Row 17 (Fel. Assault) ----
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 25
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 168
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = 51.351351351351
$ws.Range("L17").Value = 69.696969696969
$ws.Range("M17").Value = 136.619718309859
$ws.Range("N17").Value = -33.858267716535

# ---- Row 18 (Burglary) ----
# C18 switches from the text placeholder "0" to a real number.
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 142
$ws.Range("J18").Value = 129
$ws.Range("K18").Value = 10.077519379845
$ws.Range("L18").Value = 63.218390804597
$ws.Range("M18").Value = 89.333333333333
$ws.Range("N18").Value = -73.358348968105

# ---- Row 19 (Gr. Larceny) ----
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 7.272727272727
$ws.Range("I19").Value = 443
$ws.Range("J19").Value = 406
$ws.Range("K19").Value = 9.113300492610
$ws.Range("L19").Value = 41.533546325878
$ws.Range("M19").Value = 46.204620462046
$ws.Range("N19").Value = -36.350574712643

# ---- Row 20 (G.L.A.) ----
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = 38.888888888888
$ws.Range("L20").Value = 127.272727272727
$ws.Range("M20").Value = 47.058823529411
$ws.Range("N20").Value = -89.035087719298

# ---- Row 21 (TOTAL) ----
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 34.042553191489
$ws.Range("I21").Value = 901
$ws.Range("J21").Value = 784
$ws.Range("K21").Value = 14.923469387755
$ws.Range("L21").Value = 47.222222222222
$ws.Range("M21").Value = 49.419568822553
$ws.Range("N21").Value = -67.378711078928

# ---- Row 22 (Transit) ----
# C22, D22 switch from numbers to the text placeholder "0"; E22 to "***.*".
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("L22").Value = 9.677419354838

# ---- Row 23 (Housing) ----
# C23 switches from a number to the text placeholder "0".
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 150
$ws.Range("L23").Value = -3.125

# ---- Row 24 (Petit Larceny) ----
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 54.545454545454
$ws.Range("F24").Value = 170
$ws.Range("G24").Value = 170
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1319
$ws.Range("J24").Value = 1164
$ws.Range("K24").Value = 13.316151202749
$ws.Range("L24").Value = 57.398568019093
$ws.Range("M24").Value = 34.454638124362

# ---- Row 25 (Misd. Assault) ----
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 21.212121212121
$ws.Range("I25").Value = 249
$ws.Range("J25").Value = 235
$ws.Range("K25").Value = 5.957446808510
$ws.Range("L25").Value = 50
$ws.Range("M25").Value = 10.176991150442

# ---- Row 26 (UCR Rape*) ----
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = -62.5

# ---- Row 27 (Other Sex Crimes) ----
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 36
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -21.739130434782
$ws.Range("L27").Value = -10

# ---- Row 30 (Hate Crimes) ----
# G30 switches from a number to the text placeholder "0"; H30 to "***.*".
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
